{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per the supplied diff):\n//   1. Title paragraph:    \"Title\"    -> \"Title del Art\u00edculo\"\n//   2. Subtitle paragraph: \"Subtitle\" -> \"La Civilizaci\u00f3n Occidental\"\n//   3. The first \"First Paragraph\"-styled paragraph (right after the\n//      \"Heading 1\" paragraph) that begins with \"All human beings\" has\n//      that leading phrase swapped for \"Todos los seres humanos\",\n//      leaving the remainder of the paragraph untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load what we need to identify each target paragraph.\nparagraphs.items.forEach((p) => p.load(\"style,text\"));\nawait context.sync();\n\nlet titleDone = false;\nlet subtitleDone = false;\nlet firstParaDone = false;\n\nconst OLD_LEAD = \"All human beings\";\nconst NEW_LEAD = \"Todos los seres humanos\";\n\nfor (const p of paragraphs.items) {\n  if (!titleDone && p.style === \"Title\" && p.text === \"Title\") {\n    // Append \" del Art\u00edculo\" to the existing \"Title\" text.\n    p.insertText(\" del Art\u00edculo\", Word.InsertLocation.end);\n    titleDone = true;\n    continue;\n  }\n\n  if (!subtitleDone && p.style === \"Subtitle\" && p.text === \"Subtitle\") {\n    // Replace the whole \"Subtitle\" text with the new subtitle.\n    p.insertText(\"La Civilizaci\u00f3n Occidental\", Word.InsertLocation.replace);\n    subtitleDone = true;\n    continue;\n  }\n\n  if (\n    !firstParaDone &&\n    p.style === \"First Paragraph\" &&\n    p.text.indexOf(OLD_LEAD) === 0\n  ) {\n    // Only swap the leading phrase; keep the rest of the paragraph intact.\n    const newText = NEW_LEAD + p.text.slice(OLD_LEAD.length);\n    p.insertText(newText, Word.InsertLocation.replace);\n    firstParaDone = true;\n    continue;\n  }\n\n  if (titleDone && subtitleDone && firstParaDone) {\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (per the supplied diff):\n#   1. Title paragraph:    \"Title\"    -> \"Title del Art\u00edculo\"\n#   2. Subtitle paragraph: \"Subtitle\" -> \"La Civilizaci\u00f3n Occidental\"\n#   3. The first \"First Paragraph\"-styled paragraph (right after the\n#      \"Heading 1\" paragraph) that begins with \"All human beings\" has\n#      that leading phrase swapped for \"Todos los seres humanos\",\n#      leaving the remainder of the paragraph untouched.\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n$titleDone = $false\n$subtitleDone = $false\n$firstParaDone = $false\n\nfor ($i = 1; $i -le $count; $i++) {\n    if ($titleDone -and $subtitleDone -and $firstParaDone) {\n        break\n    }\n\n    $p = $d.Paragraphs($i)\n    $styleName = $p.Style.NameLocal\n    # Range.Text carries the trailing paragraph mark (chr 13) - strip it before comparing.\n    $text = $p.Range.Text.TrimEnd([char]13)\n\n    if (-not $titleDone -and $styleName -eq \"Title\" -and $text -eq \"Title\") {\n        $r = $p.Range\n        $r.End = $r.End - 1\n        $r.InsertAfter(\" del Art\u00edculo\")\n        $titleDone = $true\n        continue\n    }\n\n    if (-not $subtitleDone -and $styleName -eq \"Subtitle\" -and $text -eq \"Subtitle\") {\n        $r = $p.Range\n        $r.End = $r.End - 1\n        $r.Text = \"La Civilizaci\u00f3n Occidental\"\n        $subtitleDone = $true\n        continue\n    }\n\n    if (-not $firstParaDone -and $styleName -eq \"First Paragraph\" -and $text.StartsWith(\"All human beings\")) {\n        $r = $p.Range\n        $find = $r.Find\n        $find.Text = \"All human beings\"\n        $find.Replacement.Text = \"Todos los seres humanos\"\n        # wdFindContinue (1) keeps the search within $r; wdReplaceOne (1) replaces only the first hit.\n        $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n        $firstParaDone = $true\n        continue\n    }\n}\n"}
